# Fixing broken link on schedule page
# Delete the hidden slide (SlideID 389) that held a picture with a broken
# image link ("Minimize-interval-list solution" was slide 3, this empty
# title/content placeholder slide with the broken "people smiling" photo
# immediately followed the title slide).
$p = $ppt.ActivePresentation

$targetSlideId = 389
$slideToDelete = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    if ($candidate.SlideID -eq $targetSlideId) {
        $slideToDelete = $candidate
        break
    }
}

if ($slideToDelete -ne $null) {
    $slideToDelete.Delete()
} else {
    # Fallback: the broken-link slide is the second slide in the deck.
    $p.Slides.Item(2).Delete()
}
